$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1681.3235
$ws.Range("I15").Value = 1681.3235
$ws.Range("K15").Value = 5043.970499999999
$ws.Range("M15").Value = -4874.970499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1139
$ws.Range("I19").Value = 1458.1428
$ws.Range("J19").Value = 692.2
$ws.Range("K19").Value = 1458.1428
$ws.Range("L19").Value = 692.2
$ws.Range("M19").Value = -1283.1428
$ws.Range("N19").Value = -1042.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9755.186
$ws.Range("I69").Value = 3000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8126

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 9755.186
$ws.Range("I72").Value = 3000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -22632

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1072.8889
$ws.Range("I96").Value = 2040.25
$ws.Range("J96").Value = 299
$ws.Range("K96").Value = 6120.75
$ws.Range("L96").Value = 897
$ws.Range("M96").Value = -4747.75
$ws.Range("N96").Value = -3643

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1769.25
$ws.Range("J97").Value = 1781.4445
$ws.Range("L97").Value = 5344.333500000001
$ws.Range("N97").Value = -6336.333500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 826.3333
$ws.Range("J99").Value = 498.66666
$ws.Range("L99").Value = 1495.99998
$ws.Range("N99").Value = -4491.999980000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2542.2083
$ws.Range("I100").Value = 1883.1177
$ws.Range("K100").Value = 1883.1177
$ws.Range("M100").Value = -1342.1177

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 108189.5
$ws.Range("J117").Value = 108189.5
$ws.Range("L117").Value = 108189.5
$ws.Range("N117").Value = -117367.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 83572.25
$ws.Range("J120").Value = 83572.25
$ws.Range("L120").Value = 83572.25
$ws.Range("N120").Value = -93248.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 85537.39999999999
$ws.Range("J137").Value = 4312.3335
$ws.Range("L137").Value = 12937.0005
$ws.Range("N137").Value = -18037.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3897.3547
$ws.Range("I61").Value = 3258.5715
$ws.Range("K61").Value = 3258.5715
$ws.Range("M61").Value = -3046.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1707.3043
$ws.Range("I97").Value = 1221.7059
$ws.Range("K97").Value = 1221.7059
$ws.Range("M97").Value = -725.7058999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 36795.332
$ws.Range("J112").Value = 36795.332
$ws.Range("L112").Value = 36795.332
$ws.Range("N112").Value = -39749.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 280728.6
$ws.Range("I132").Value = 305943.6
$ws.Range("K132").Value = 917830.7999999999
$ws.Range("M132").Value = -915300.7999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3897.3547
$ws.Range("I136").Value = 3258.5715
$ws.Range("K136").Value = 9775.7145
$ws.Range("M136").Value = -7225.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 59346.668
$ws.Range("J51").Value = 59346.668
$ws.Range("L51").Value = 59346.668
$ws.Range("N51").Value = -60328.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 765.8461
$ws.Range("J86").Value = 744.8570999999999
$ws.Range("L86").Value = 744.8570999999999
$ws.Range("N86").Value = -2990.8571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 765.8461
$ws.Range("J89").Value = 744.8570999999999
$ws.Range("L89").Value = 3724.2855
$ws.Range("N89").Value = -14956.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 584.8461
$ws.Range("I94").Value = 564.8570999999999
$ws.Range("J94").Value = 608.1667
$ws.Range("K94").Value = 564.8570999999999
$ws.Range("L94").Value = 608.1667
$ws.Range("M94").Value = -113.8570999999999
$ws.Range("N94").Value = -1510.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3907.4546
$ws.Range("I105").Value = 4666.3335
$ws.Range("K105").Value = 4666.3335
$ws.Range("M105").Value = -2919.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 125234.5
$ws.Range("J112").Value = 125234.5
$ws.Range("L112").Value = 125234.5
$ws.Range("N112").Value = -128188.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 115013.5
$ws.Range("J120").Value = 115013.5
$ws.Range("L120").Value = 115013.5
$ws.Range("N120").Value = -124689.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2611.7778
$ws.Range("I62").Value = 1958.5714
$ws.Range("K62").Value = 1958.5714
$ws.Range("M62").Value = -1334.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2611.7778
$ws.Range("I65").Value = 1958.5714
$ws.Range("K65").Value = 9792.857
$ws.Range("M65").Value = -6672.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2679.5715
$ws.Range("J99").Value = 2927.8572
$ws.Range("L99").Value = 2927.8572
$ws.Range("N99").Value = -5923.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3111.3103
$ws.Range("J122").Value = 3500.4666
$ws.Range("L122").Value = 10501.3998
$ws.Range("N122").Value = -15401.3998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2679.5715
$ws.Range("J126").Value = 2927.8572
$ws.Range("L126").Value = 8783.571599999999
$ws.Range("N126").Value = -13723.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2666
$ws.Range("I132").Value = 2677.0667
$ws.Range("K132").Value = 8031.2001
$ws.Range("M132").Value = -5501.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1886.8889
$ws.Range("I134").Value = 1872.875
$ws.Range("K134").Value = 5618.625
$ws.Range("M134").Value = -3083.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6306.174
$ws.Range("J137").Value = 11189.818
$ws.Range("L137").Value = 33569.454
$ws.Range("N137").Value = -43769.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2167.6
$ws.Range("J80").Value = 1794.2858
$ws.Range("L80").Value = 1794.2858
$ws.Range("N80").Value = -3790.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2167.6
$ws.Range("J83").Value = 1794.2858
$ws.Range("L83").Value = 8971.429
$ws.Range("N83").Value = -18955.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 96997.664
$ws.Range("J111").Value = 96997.664
$ws.Range("L111").Value = 96997.664
$ws.Range("N111").Value = -103131.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2038.8
$ws.Range("I40").Value = 2038.8
$ws.Range("K40").Value = 2038.8
$ws.Range("M40").Value = -1902.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3782
$ws.Range("J46").Value = 3815.9565
$ws.Range("L46").Value = 3815.9565
$ws.Range("N46").Value = -4191.9565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 74146.39999999999
$ws.Range("J110").Value = 74146.39999999999
$ws.Range("L110").Value = 74146.39999999999
$ws.Range("N110").Value = -82326.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 17204.572
$ws.Range("I122").Value = 17994.666
$ws.Range("J122").Value = 15782.4
$ws.Range("K122").Value = 53983.99800000001
$ws.Range("L122").Value = 47347.2
$ws.Range("M122").Value = -51533.99800000001
$ws.Range("N122").Value = -52247.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 456129.8
$ws.Range("I132").Value = 527845.0600000001
$ws.Range("K132").Value = 1583535.18
$ws.Range("M132").Value = -1581005.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 4000
$ws.Range("K39").Value = 4000
$ws.Range("M39").Value = -3587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6064
$ws.Range("I62").Value = 4166.6665
$ws.Range("J62").Value = 6877.143
$ws.Range("K62").Value = 4166.6665
$ws.Range("L62").Value = 6877.143
$ws.Range("N62").Value = -8125.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6064
$ws.Range("I65").Value = 4166.6665
$ws.Range("J65").Value = 6877.143
$ws.Range("K65").Value = 20833.3325
$ws.Range("L65").Value = 34385.715
$ws.Range("N65").Value = -40625.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35382.03
$ws.Range("I132").Value = 41155.73
$ws.Range("K132").Value = 123467.19
$ws.Range("M132").Value = -120937.19

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2290.8386
$ws.Range("I136").Value = 1433.2632
$ws.Range("K136").Value = 4299.7896
$ws.Range("M136").Value = -1749.7896
